$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 7.240679918064298
$ws.Cells.Item(2, 4).Value = 5.82196662681687
$ws.Cells.Item(2, 5).Value = 12.46065815797078
$ws.Cells.Item(2, 6).Value = 53.88144777410852
$ws.Cells.Item(2, 7).Value = 71.55053117575966
$ws.Cells.Item(2, 8).Value = 24.17023702721616
$ws.Cells.Item(2, 9).Value = 40.98355024799
$ws.Cells.Item(2, 10).Value = 10.92775514833054
$ws.Cells.Item(2, 12).Value = 10.45052131975543
$ws.Cells.Item(2, 14).Value = 18.99769123883691
$ws.Cells.Item(3, 3).Value = 7.253979846253294
$ws.Cells.Item(3, 4).Value = 5.803469237481298
$ws.Cells.Item(3, 5).Value = 12.4895881847224
$ws.Cells.Item(3, 6).Value = 53.47938016116156
$ws.Cells.Item(3, 7).Value = 70.64340532633348
$ws.Cells.Item(3, 8).Value = 24.08215247237939
$ws.Cells.Item(3, 9).Value = 40.73299499319394
$ws.Cells.Item(3, 10).Value = 10.95296466775718
$ws.Cells.Item(3, 12).Value = 10.47620201549809
$ws.Cells.Item(3, 14).Value = 18.40031508502705
$ws.Cells.Item(4, 3).Value = 7.262716213273668
$ws.Cells.Item(4, 4).Value = 5.791878923516001
$ws.Cells.Item(4, 5).Value = 12.50899615333653
$ws.Cells.Item(4, 6).Value = 53.24722460944285
$ws.Cells.Item(4, 7).Value = 70.10457091085885
$ws.Cells.Item(4, 8).Value = 24.03432458196443
$ws.Cells.Item(4, 9).Value = 40.59028273905867
$ws.Cells.Item(4, 10).Value = 10.97000095944372
$ws.Cells.Item(4, 12).Value = 10.49323879232871
$ws.Cells.Item(4, 14).Value = 18.02485520896364
$ws.Cells.Item(5, 3).Value = 7.266420037099855
$ws.Cells.Item(5, 4).Value = 5.787096938839175
$ws.Cells.Item(5, 5).Value = 12.51731877923831
$ws.Cells.Item(5, 6).Value = 53.15639101174041
$ws.Cells.Item(5, 7).Value = 69.88979535565943
$ws.Cells.Item(5, 8).Value = 24.01641548997424
$ws.Cells.Item(5, 9).Value = 40.5349625505344
$ws.Cells.Item(5, 10).Value = 10.97733491058785
$ws.Cells.Item(5, 12).Value = 10.50050077123841
$ws.Cells.Item(5, 14).Value = 17.86990355188765
$ws.Cells.Item(6, 3).Value = 7.267043741598322
$ws.Cells.Item(6, 4).Value = 5.78629933179792
$ws.Cells.Item(6, 5).Value = 12.51872573618084
$ws.Cells.Item(6, 6).Value = 53.14153789802702
$ws.Cells.Item(6, 7).Value = 69.85442899561212
$ws.Cells.Item(6, 8).Value = 24.01353738982078
$ws.Cells.Item(6, 9).Value = 40.52594897505973
$ws.Cells.Item(6, 10).Value = 10.97857635071194
$ws.Cells.Item(6, 12).Value = 10.50172591506596
$ws.Cells.Item(6, 14).Value = 17.84406337566597
$ws.Cells.Item(7, 3).Value = 7.262765582145795
$ws.Cells.Item(7, 4).Value = 5.791814670616975
$ws.Cells.Item(7, 5).Value = 12.50910671986518
$ws.Cells.Item(7, 6).Value = 53.24598423564453
$ws.Cells.Item(7, 7).Value = 70.10165461186129
$ws.Cells.Item(7, 8).Value = 24.03407664256463
$ws.Cells.Item(7, 9).Value = 40.58952514192086
$ws.Cells.Item(7, 10).Value = 10.97009828257184
$ws.Cells.Item(7, 12).Value = 10.49333543642013
$ws.Cells.Item(7, 14).Value = 18.02277304766463
$ws.Cells.Item(8, 3).Value = 7.245147584407842
$ws.Cells.Item(8, 4).Value = 5.815636878004722
$ws.Cells.Item(8, 5).Value = 12.47029192098302
$ws.Cells.Item(8, 6).Value = 53.73980567169919
$ws.Cells.Item(8, 7).Value = 71.2341364771066
$ws.Cells.Item(8, 8).Value = 24.13857231598888
$ws.Cells.Item(8, 9).Value = 40.89487285550817
$ws.Cells.Item(8, 10).Value = 10.93612404254448
$ws.Cells.Item(8, 12).Value = 10.4591128923743
$ws.Cells.Item(8, 14).Value = 18.79364780656866
$ws.Cells.Item(9, 3).Value = 7.215108822380648
$ws.Cells.Item(9, 4).Value = 5.860518890367791
$ws.Cells.Item(9, 5).Value = 12.40722232928409
$ws.Cells.Item(9, 6).Value = 54.82147432476117
$ws.Cells.Item(9, 7).Value = 73.58779030816854
$ws.Cells.Item(9, 8).Value = 24.39269338200491
$ws.Cells.Item(9, 9).Value = 41.58009133903944
$ws.Cells.Item(9, 10).Value = 10.88186553487712
$ws.Cells.Item(9, 12).Value = 10.40205533037962
$ws.Cells.Item(9, 14).Value = 20.2273683202997
$ws.Cells.Item(10, 3).Value = 7.195769984866931
$ws.Cells.Item(10, 4).Value = 5.8923976990045
$ws.Cells.Item(10, 5).Value = 12.36883266751566
$ws.Cells.Item(10, 6).Value = 55.68014242145632
$ws.Cells.Item(10, 7).Value = 75.38261594773743
$ws.Cells.Item(10, 8).Value = 24.60868663556973
$ws.Cells.Item(10, 9).Value = 42.13339779608805
$ws.Cells.Item(10, 10).Value = 10.84954933826331
$ws.Cells.Item(10, 12).Value = 10.36624493908714
$ws.Cells.Item(10, 14).Value = 21.22223697909767
$ws.Cells.Item(11, 3).Value = 7.187561265708002
$ws.Cells.Item(11, 4).Value = 5.906666910723553
$ws.Cells.Item(11, 5).Value = 12.3530928670554
$ws.Cells.Item(11, 6).Value = 56.0834620339088
$ws.Cells.Item(11, 7).Value = 76.21001617407893
$ws.Cells.Item(11, 8).Value = 24.71311884555069
$ws.Cells.Item(11, 9).Value = 42.39525773261941
$ws.Cells.Item(11, 10).Value = 10.83648883126462
$ws.Cells.Item(11, 12).Value = 10.3512767649242
$ws.Cells.Item(11, 14).Value = 21.66018057919901
$ws.Cells.Item(12, 3).Value = 7.184537182689872
$ws.Cells.Item(12, 4).Value = 5.912037162253204
$ws.Cells.Item(12, 5).Value = 12.34738042714542
$ws.Cells.Item(12, 6).Value = 56.23790833005934
$ws.Cells.Item(12, 7).Value = 76.5246205168887
$ws.Cells.Item(12, 8).Value = 24.75353391884827
$ws.Cells.Item(12, 9).Value = 42.49581509440111
$ws.Cells.Item(12, 10).Value = 10.8317792004487
$ws.Cells.Item(12, 12).Value = 10.34579855375347
$ws.Cells.Item(12, 14).Value = 21.82377585682186
$ws.Cells.Item(13, 3).Value = 7.185184724326807
$ws.Cells.Item(13, 4).Value = 5.910882054456976
$ws.Cells.Item(13, 5).Value = 12.34859967880871
$ws.Cells.Item(13, 6).Value = 56.20457092824671
$ws.Cells.Item(13, 7).Value = 76.45681216362989
$ws.Cells.Item(13, 8).Value = 24.74479148950045
$ws.Cells.Item(13, 9).Value = 42.47409719890825
$ws.Cells.Item(13, 10).Value = 10.83278299952341
$ws.Cells.Item(13, 12).Value = 10.34696994147886
$ws.Cells.Item(13, 14).Value = 21.78864458690801
$ws.Cells.Item(14, 3).Value = 7.187310782939388
$ws.Cells.Item(14, 4).Value = 5.907109388213589
$ws.Cells.Item(14, 5).Value = 12.35261793331732
$ws.Cells.Item(14, 6).Value = 56.0961346323378
$ws.Cells.Item(14, 7).Value = 76.2358747157632
$ws.Cells.Item(14, 8).Value = 24.71642652038864
$ws.Cells.Item(14, 9).Value = 42.40350303253088
$ws.Cells.Item(14, 10).Value = 10.83609663402414
$ws.Cells.Item(14, 12).Value = 10.35082226427555
$ws.Cells.Item(14, 14).Value = 21.67368539489659
$ws.Cells.Item(15, 3).Value = 7.188624036432071
$ws.Cells.Item(15, 4).Value = 5.904794206737772
$ws.Cells.Item(15, 5).Value = 12.35511151146535
$ws.Cells.Item(15, 6).Value = 56.02993468999172
$ws.Cells.Item(15, 7).Value = 76.10070318765389
$ws.Cells.Item(15, 8).Value = 24.69916471176785
$ws.Cells.Item(15, 9).Value = 42.36044200767237
$ws.Cells.Item(15, 10).Value = 10.83815708467463
$ws.Cells.Item(15, 12).Value = 10.35320664912925
$ws.Cells.Item(15, 14).Value = 21.60297336126124
$ws.Cells.Item(16, 3).Value = 7.19631826594784
$ws.Cells.Item(16, 4).Value = 5.891460513898252
$ws.Cells.Item(16, 5).Value = 12.36989598816727
$ws.Cells.Item(16, 6).Value = 55.65403096743166
$ws.Cells.Item(16, 7).Value = 75.32873880181515
$ws.Cells.Item(16, 8).Value = 24.6019844126132
$ws.Cells.Item(16, 9).Value = 42.11648363326397
$ws.Cells.Item(16, 10).Value = 10.85043589954232
$ws.Cells.Item(16, 12).Value = 10.36724973123358
$ws.Cells.Item(16, 14).Value = 21.19330956972086
$ws.Cells.Item(17, 3).Value = 7.201188997336454
$ws.Cells.Item(17, 4).Value = 5.88322133187757
$ws.Cells.Item(17, 5).Value = 12.37940729992349
$ws.Cells.Item(17, 6).Value = 55.42660414023256
$ws.Cells.Item(17, 7).Value = 74.857758199261
$ws.Cells.Item(17, 8).Value = 24.54393606585015
$ws.Cells.Item(17, 9).Value = 41.96938018930108
$ws.Cells.Item(17, 10).Value = 10.85838883625898
$ws.Cells.Item(17, 12).Value = 10.37620318873222
$ws.Cells.Item(17, 14).Value = 20.93814219015166
$ws.Cells.Item(18, 3).Value = 7.204045931744877
$ws.Cells.Item(18, 4).Value = 5.878460461598528
$ws.Cells.Item(18, 5).Value = 12.38504021210554
$ws.Cells.Item(18, 6).Value = 55.2969971001227
$ws.Cells.Item(18, 7).Value = 74.58790791313899
$ws.Cells.Item(18, 8).Value = 24.51113074409971
$ws.Cells.Item(18, 9).Value = 41.88573033908828
$ws.Cells.Item(18, 10).Value = 10.86311753393232
$ws.Cells.Item(18, 12).Value = 10.38147744019941
$ws.Cells.Item(18, 14).Value = 20.79000725568362
$ws.Cells.Item(19, 3).Value = 7.205022766883077
$ws.Cells.Item(19, 4).Value = 5.876844741606675
$ws.Cells.Item(19, 5).Value = 12.3869752874135
$ws.Cells.Item(19, 6).Value = 55.25332424940351
$ws.Cells.Item(19, 7).Value = 74.49672950215897
$ws.Cells.Item(19, 8).Value = 24.50012403820437
$ws.Cells.Item(19, 9).Value = 41.85757474748316
$ws.Cells.Item(19, 10).Value = 10.86474509740382
$ws.Cells.Item(19, 12).Value = 10.38328459373762
$ws.Cells.Item(19, 14).Value = 20.73962067985785
$ws.Cells.Item(20, 3).Value = 7.200664765968112
$ws.Cells.Item(20, 4).Value = 5.884100672650497
$ws.Cells.Item(20, 5).Value = 12.37837801159877
$ws.Cells.Item(20, 6).Value = 55.45069040463455
$ws.Cells.Item(20, 7).Value = 74.90778877801738
$ws.Cells.Item(20, 8).Value = 24.55005525359931
$ws.Cells.Item(20, 9).Value = 41.98494068367541
$ws.Cells.Item(20, 10).Value = 10.85752625329142
$ws.Cells.Item(20, 12).Value = 10.37523719837484
$ws.Cells.Item(20, 14).Value = 20.96544799483449
$ws.Cells.Item(21, 3).Value = 7.186684019975634
$ws.Cells.Item(21, 4).Value = 5.908218410409182
$ws.Cells.Item(21, 5).Value = 12.35143094719327
$ws.Cells.Item(21, 6).Value = 56.12793926527856
$ws.Cells.Item(21, 7).Value = 76.30073676084855
$ws.Cells.Item(21, 8).Value = 24.7247345797031
$ws.Cells.Item(21, 9).Value = 42.42420086915446
$ws.Cells.Item(21, 10).Value = 10.83511692977273
$ws.Cells.Item(21, 12).Value = 10.34968559119871
$ws.Cells.Item(21, 14).Value = 21.70751365554064
$ws.Cells.Item(22, 3).Value = 7.178038512520449
$ws.Cells.Item(22, 4).Value = 5.923787415286124
$ws.Cells.Item(22, 5).Value = 12.33526432029697
$ws.Cells.Item(22, 6).Value = 56.58052096965619
$ws.Cells.Item(22, 7).Value = 77.21847383996338
$ws.Cells.Item(22, 8).Value = 24.84395200054346
$ws.Cells.Item(22, 9).Value = 42.71939213923282
$ws.Cells.Item(22, 10).Value = 10.8218474141848
$ws.Cells.Item(22, 12).Value = 10.33409299526921
$ws.Cells.Item(22, 14).Value = 22.17935961385674
$ws.Cells.Item(23, 3).Value = 7.182607875268259
$ws.Cells.Item(23, 4).Value = 5.915495547875958
$ws.Cells.Item(23, 5).Value = 12.34376055245106
$ws.Cells.Item(23, 6).Value = 56.33809516110047
$ws.Cells.Item(23, 7).Value = 76.72807974620693
$ws.Cells.Item(23, 8).Value = 24.77986769463015
$ws.Cells.Item(23, 9).Value = 42.56112286494051
$ws.Cells.Item(23, 10).Value = 10.82880360110862
$ws.Cells.Item(23, 12).Value = 10.34231384795075
$ws.Cells.Item(23, 14).Value = 21.92877110911181
$ws.Cells.Item(24, 3).Value = 7.200901594385742
$ws.Cells.Item(24, 4).Value = 5.883703197632071
$ws.Cells.Item(24, 5).Value = 12.37884283974535
$ws.Cells.Item(24, 6).Value = 55.43979743396029
$ws.Cells.Item(24, 7).Value = 74.88516705758251
$ws.Cells.Item(24, 8).Value = 24.54728699991263
$ws.Cells.Item(24, 9).Value = 41.97790290789958
$ws.Cells.Item(24, 10).Value = 10.85791573976648
$ws.Cells.Item(24, 12).Value = 10.3756735276783
$ws.Cells.Item(24, 14).Value = 20.95310750188673
$ws.Cells.Item(25, 3).Value = 7.222754232826133
$ws.Cells.Item(25, 4).Value = 5.848571135465709
$ws.Cells.Item(25, 5).Value = 12.42288824186712
$ws.Cells.Item(25, 6).Value = 54.5172289374464
$ws.Cells.Item(25, 7).Value = 72.93841039704067
$ws.Cells.Item(25, 8).Value = 24.31874453786596
$ws.Cells.Item(25, 9).Value = 41.38574222493524
$ws.Cells.Item(25, 10).Value = 10.8952191525488
$ws.Cells.Item(25, 12).Value = 10.41641664866672
$ws.Cells.Item(25, 14).Value = 19.84905939529497
